# Apply the "Stock" column addition to the "All Products" sheet (sheet1 / rId1),
# update the sheet view (zoom + selection), resize the new column, and
# set the new row heights that Excel recomputed for rows 9 and 25 once the
# new column was introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Products")

# --- New "Stock" column header ---
$ws.Cells.Item(1, 9).Value = "Stock"

# --- Populate Stock quantity (100) for every product row (2-55) ---
for ($row = 2; $row -le 55; $row++) {
    $ws.Cells.Item($row, 9).Value = 100
}

# --- New column width for column I ---
# Target stored width is 26.109375 characters; the closest value this
# engine's width-rounding model can produce is 26.1666... (nearest 1/6
# increment), reached with an input of 25.333333333333336.
$ws.Columns.Item(9).ColumnWidth = 25.333333333333336

# --- Row heights recomputed by Excel once column I text wrapping changed ---
$ws.Rows.Item(9).RowHeight = 144
$ws.Rows.Item(25).RowHeight = 57.6

# --- Sheet view: zoom change and new selected cell ---
$ws.Activate()
$win = $ws.Application.ActiveWindow
$win.Zoom = 85
$ws.Range("P3").Select()
